$wb = $excel.ActiveWorkbook

# New build/version string for this release (was: "mines - January 30 (built on February 02 2026 12.49.33 EST)")
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Update "Version:" line on the About sheet (A2)
$wsAbout.Range("A2").Value = "Version: " + $newVersion

# Update the Recommended Citation line on the About sheet (A6)
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Beckley Pocahontas Coal Mine, United States, M0990, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Update the build_version column (S) for every data row on the data sheet
$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
